# Auto-generated Excel COM-interop script applying the Alexander_Profits diff.
# For each changed cell we set the new value directly; for cells that are
# removed entirely in the diff (blank after edit) we clear them so the
# resulting OOXML omits the <c> element, matching openpyxl/Excel behaviour
# for blank cells.

$wb = $excel.ActiveWorkbook


# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H98").Value = 1568.1428
$ws.Range("I98").Value = 889.3333
$ws.Range("J98").Value = 2790
$ws.Range("K98").Value = 889.3333
$ws.Range("L98").Value = 2790
$ws.Range("M98").Value = 608.6667
$ws.Range("N98").Value = -5786
$ws.Range("H122").Value = 1568.1428
$ws.Range("I122").Value = 889.3333
$ws.Range("J122").Value = 2790
$ws.Range("K122").Value = 2667.9999
$ws.Range("L122").Value = 8370
$ws.Range("M122").Value = -217.9998999999998
$ws.Range("N122").Value = -13270

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 1017.5
$ws.Range("I74").Value = 1026.6428
$ws.Range("K74").Value = 1026.6428
$ws.Range("M74").Value = -152.6428000000001
$ws.Range("H77").Value = 1017.5
$ws.Range("I77").Value = 1026.6428
$ws.Range("K77").Value = 5133.214
$ws.Range("M77").Value = -765.2139999999999
$ws.Range("H102").Value = 83335464
$ws.Range("I102").Value = 2289.889
$ws.Range("J102").Value = 333334980
$ws.Range("K102").Value = 2289.889
$ws.Range("L102").Value = 333334980
$ws.Range("M102").Value = -667.8890000000001
$ws.Range("N102").Value = -333338224
$ws.Range("H122").Value = 2117.3333
$ws.Range("I122").Value = 2065.4119
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6196.2357
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3746.2357
$ws.Range("N122").Value = -13900

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 3532.6667
$ws.Range("I105").Value = 1969
$ws.Range("J105").Value = 6660
$ws.Range("K105").Value = 1969
$ws.Range("L105").Value = 6660
$ws.Range("M105").Value = -222
$ws.Range("N105").Value = -10154

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H36").Value = 20309.6
$ws.Range("I36").Value = 3548
$ws.Range("J36").Value = 24500
$ws.Range("K36").Value = 3548
$ws.Range("L36").Value = 24500
$ws.Range("M36").Value = -3160
$ws.Range("N36").Value = -25276
$ws.Range("H38").Value = 48500
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""   # was -1598 -> blank
$ws.Range("H39").Value = 32333.334
$ws.Range("I39").Value = 32333.334
$ws.Range("K39").Value = 32333.334
$ws.Range("M39").Value = -31942.334
$ws.Range("H40").Value = 20309.6
$ws.Range("I40").Value = 3548
$ws.Range("J40").Value = 24500
$ws.Range("K40").Value = 3548
$ws.Range("L40").Value = 24500
$ws.Range("M40").Value = -3388
$ws.Range("N40").Value = -24820
$ws.Range("H41").Value = 14988.333
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 19982.5
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 19982.5
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -20838.5
$ws.Range("H42").Value = 46000
$ws.Range("J42").Value = 46000
$ws.Range("L42").Value = 46000
$ws.Range("N42").Value = -47186
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = ""   # was 27777.334 -> blank
$ws.Range("M44").Value = 30000
$ws.Range("N44").Value = -30884
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = ""   # was 8500 -> blank
$ws.Range("N45").Value = 0
$ws.Range("H46").Value = 48500
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""   # was -1764 -> blank
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = ""   # was 6000 -> blank
$ws.Range("N47").Value = 0
$ws.Range("H48").Value = 5250
$ws.Range("I48").Value = 4500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -4024
$ws.Range("H49").Value = 32333.334
$ws.Range("I49").Value = 32333.334
$ws.Range("K49").Value = 32333.334
$ws.Range("M49").Value = -32151.334
$ws.Range("H50").Value = 9312.6
$ws.Range("J50").Value = 9312.6
$ws.Range("L50").Value = 9312.6
$ws.Range("N50").Value = -10562.6
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = ""   # was 9342.833000000001 -> blank
$ws.Range("N51").Value = 0
$ws.Range("H54").Value = 5452
$ws.Range("J54").Value = 5452
$ws.Range("L54").Value = 5452
$ws.Range("N54").Value = -6768
$ws.Range("H55").Value = 11899.833
$ws.Range("I55").Value = 5999.5
$ws.Range("J55").Value = 14850
$ws.Range("K55").Value = 5999.5
$ws.Range("L55").Value = 14850
$ws.Range("M55").Value = -5684.5
$ws.Range("N55").Value = -15480
$ws.Range("H56").Value = 7697.5
$ws.Range("J56").Value = 7697.5
$ws.Range("L56").Value = 7697.5
$ws.Range("N56").Value = -9387.5
$ws.Range("H57").Value = 7993.3335
$ws.Range("J57").Value = 7993.3335
$ws.Range("L57").Value = 7993.3335
$ws.Range("N57").Value = -9113.333500000001
$ws.Range("H59").Value = 16027
$ws.Range("J59").Value = 16027
$ws.Range("L59").Value = 16027
$ws.Range("N59").Value = -18317
$ws.Range("H60").Value = 276823.75
$ws.Range("J60").Value = 276823.75
$ws.Range("L60").Value = 276823.75
$ws.Range("N60").Value = -277845.75
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = ""   # was 9342.833000000001 -> blank
$ws.Range("N61").Value = 0
$ws.Range("H63").Value = 33271
$ws.Range("J63").Value = 33271
$ws.Range("L63").Value = 33271
$ws.Range("N63").Value = -34643
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = ""   # was 33000 -> blank
$ws.Range("N64").Value = 0
$ws.Range("H66").Value = 33271
$ws.Range("J66").Value = 33271
$ws.Range("L66").Value = 99813
$ws.Range("N66").Value = -106677
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = ""   # was 33000 -> blank
$ws.Range("N67").Value = 0
$ws.Range("H68").Value = 17999
$ws.Range("J68").Value = 17999
$ws.Range("L68").Value = 17999
$ws.Range("N68").Value = -19497
$ws.Range("H70").Value = 32000
$ws.Range("J70").Value = 32000
$ws.Range("L70").Value = 32000
$ws.Range("N70").Value = -32630
$ws.Range("H71").Value = 17999
$ws.Range("J71").Value = 17999
$ws.Range("L71").Value = 53997
$ws.Range("N71").Value = -61485
$ws.Range("H73").Value = 32000
$ws.Range("J73").Value = 32000
$ws.Range("L73").Value = 32000
$ws.Range("N73").Value = -34184

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 25209.863
$ws.Range("I70").Value = 27904.732
$ws.Range("J70").Value = 4998.3335
$ws.Range("K70").Value = 27904.732
$ws.Range("L70").Value = 4998.3335
$ws.Range("M70").Value = -27634.732
$ws.Range("N70").Value = -5538.3335
$ws.Range("H73").Value = 25209.863
$ws.Range("I73").Value = 27904.732
$ws.Range("J73").Value = 4998.3335
$ws.Range("K73").Value = 27904.732
$ws.Range("L73").Value = 4998.3335
$ws.Range("M73").Value = -26968.732
$ws.Range("N73").Value = -6870.3335
$ws.Range("H122").Value = 1872.84
$ws.Range("I122").Value = 1822.6875
$ws.Range("J122").Value = 1962
$ws.Range("K122").Value = 5468.0625
$ws.Range("L122").Value = 5886
$ws.Range("M122").Value = -3018.0625
$ws.Range("N122").Value = -10786

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 2436.0527
$ws.Range("I136").Value = 2019
$ws.Range("K136").Value = 6057
$ws.Range("M136").Value = -3507
